$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update recalculated statistics values (re-write of the experiment analysis) ---
$ws.Range("C7").Value = 0.20432649999999999
$ws.Range("D7").Value = 0.02421831

$ws.Range("C8").Value = 0.2073287
$ws.Range("D8").Value = 0.02544927

$ws.Range("C9").Value = 0.19822629999999999
$ws.Range("D9").Value = 0.02511196

# --- Update the view state: scroll back to top-left and move the active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K20").Select()
